$d = $word.ActiveDocument
$t = $d.Tables(1)

# Summary rows near the top of the table were updated with new values
# (and some placeholder "0M" markers), row indices are 1-based.
$t.Rows(1).Cells(1).Range.Text  = "0M"
$t.Rows(2).Cells(1).Range.Text  = "0M"
$t.Rows(3).Cells(1).Range.Text  = "0M"
$t.Rows(4).Cells(1).Range.Text  = "1300"
$t.Rows(6).Cells(1).Range.Text  = "0.21836"
$t.Rows(7).Cells(1).Range.Text  = "0.02765"
$t.Rows(8).Cells(1).Range.Text  = "0.00558"
$t.Rows(9).Cells(1).Range.Text  = "0.07573"
$t.Rows(10).Cells(1).Range.Text = "0.07573"
$t.Rows(11).Cells(1).Range.Text = "0.10855"
$t.Rows(12).Cells(1).Range.Text = "1.82006"

# The three rows that used to hold a full tab-separated breakdown are
# collapsed down to a single summary value each.
$t.Rows(44).Cells(1).Range.Text = "99.7"
$t.Rows(45).Cells(1).Range.Text = "1.82"
$t.Rows(46).Cells(1).Range.Text = "616"
